# One-click update from Van Paper 08:46 AM on 2025-12-23
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in two previously-blank "Last Invoice Date" cells ---
# Copy number formatting from an already-formatted date cell (D2) so the
# new values inherit the same style (numFmt 165, left/top aligned) rather
# than creating a brand-new style entry.
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value2 = 46000

$ws.Range("D2").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value2 = 46009

# --- Row 27: renamed prospect + changed Prospect code ---
$ws.Range("A27").Value2 = "NATURE PATHWAYS EARLY LEARNING"
$ws.Range("C27").Value2 = "023"

# --- Row 35: Salesperson reassigned ---
$ws.Range("B35").Value2 = "Monroe, Michael D"

# --- Row 37: now a different prospect entirely ---
$ws.Range("A37").Value2 = "ENDOCRINOLOGY CLINIC OF MPLS"
$ws.Range("B37").Value2 = "Steiner, Owen A"
$ws.Range("C37").Value2 = "040"
$ws.Range("E37").Value2 = "0008396"

# --- Row 38: now a different customer (with an invoice date) ---
$ws.Range("A38").Value2 = "WOODLANE FLOWERS INC"
$ws.Range("B38").Value2 = "House Account"
$ws.Range("C38").Value2 = "013"
$ws.Range("D38").Value2 = 46003
$ws.Range("E38").Value2 = "0004759"

# --- Row 39: now a different customer (with an invoice date) ---
$ws.Range("A39").Value2 = "PILGRIM DRY CLEANERS INC"
$ws.Range("B39").Value2 = "Steiner, Owen A"
$ws.Range("C39").Value2 = "003"
$ws.Range("D39").Value2 = 46006
$ws.Range("E39").Value2 = "0004938"

# --- Append two new rows (40, 41) with the records that were pushed down ---
# Row 40 (copy formatting from row 38 which already has a populated D cell)
$ws.Range("A38:F38").Copy()
$ws.Range("A40:F40").PasteSpecial(-4122)
$ws.Rows.Item(40).RowHeight = 13.05
$ws.Range("A40").Value2 = "HOLY FAMILY MARONITE CHURCH"
$ws.Range("B40").Value2 = "Bloch, Lea L"
$ws.Range("C40").Value2 = "003"
$ws.Range("D40").Value2 = 45932
$ws.Range("E40").Value2 = "0004965"

# Row 41
$ws.Range("A38:F38").Copy()
$ws.Range("A41:F41").PasteSpecial(-4122)
$ws.Rows.Item(41).RowHeight = 13.05
$ws.Range("A41").Value2 = "SCHMITT MUSIC CTR"
$ws.Range("B41").Value2 = "Monroe, Michael D"
$ws.Range("C41").Value2 = "003"
$ws.Range("D41").Value2 = 45954
$ws.Range("E41").Value2 = "0005169"
